$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: the previously-empty E10 cell is dropped when the row is re-saved
# (E10 has no value before or after this edit - nothing visible changes).
$ws.Range("E10").Value = ""

# New row 11: a duplicate of row 10, but with "Unidades Estructura/Paneles"
# (E) filled in instead of "Paneles" (D), and "Pajareras" (M) set to "Sí"
# instead of "1".
$ws.Range("A11").Value = 2488
$ws.Range("B11").Value = "Test Ringover (NO TOCAR)"
$ws.Range("C11").Value = "Estructura coplanar NOVOTEGRA"
$ws.Range("D11").Value = ""
# Leading apostrophe forces these numeric-looking values to be stored as
# text (matching the source row, where they are plain strings, not numbers).
$ws.Range("E11").Value = "'1"
$ws.Range("F11").Value = "HUAWEI Optimizador 600W"
$ws.Range("G11").Value = "'2"
$ws.Range("H11").Value = "Inversor híbrido monofásico SUN-6k-SG05LP1-EU"
$ws.Range("I11").Value = "'1"
$ws.Range("J11").Value = "BATERÍA LITIO SIGEN ENERGY SIGENSTOR 10,0KW"
$ws.Range("K11").Value = "'3"
$ws.Range("L11").Value = "RAEDIAN CARGADOR NEO 7KW SILVER"
$ws.Range("M11").Value = "Sí"
$ws.Range("N11").Value = "2024-01-03T10:49:29.104Z"
